$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 894.5
$ws.Range("I28").Value = 543.3333
$ws.Range("J28").Value = 1421.25
$ws.Range("K28").Value = 543.3333
$ws.Range("L28").Value = 1421.25
$ws.Range("M28").Value = -58.33330000000001
$ws.Range("N28").Value = -2391.25
$ws.Range("H62").Value = 1780.8334
$ws.Range("I62").Value = 1780.8334
$ws.Range("K62").Value = 1780.8334
$ws.Range("M62").Value = -1156.8334
$ws.Range("H65").Value = 1780.8334
$ws.Range("I65").Value = 1780.8334
$ws.Range("K65").Value = 8904.166999999999
$ws.Range("M65").Value = -5784.166999999999
$ws.Range("H100").Value = 3698.1333
$ws.Range("I100").Value = 3400
$ws.Range("K100").Value = 3400
$ws.Range("M100").Value = -2859
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
$ws.Range("H132").Value = 1651.25
$ws.Range("I132").Value = 1651.25
$ws.Range("K132").Value = 4953.75
$ws.Range("M132").Value = -2423.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2732.8333
$ws.Range("I45").Value = 2732.8333
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2732.8333
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2355.8333
$ws.Range("N45").Value = $null
$ws.Range("H61").Value = 3143.75
$ws.Range("I61").Value = 3143.75
$ws.Range("K61").Value = 3143.75
$ws.Range("M61").Value = -2931.75
$ws.Range("H97").Value = 27779620
$ws.Range("I97").Value = 41668932
$ws.Range("K97").Value = 41668932
$ws.Range("M97").Value = -41668436
$ws.Range("H102").Value = 8442136
$ws.Range("I102").Value = 612916.4399999999
$ws.Range("J102").Value = 28574414
$ws.Range("K102").Value = 612916.4399999999
$ws.Range("L102").Value = 28574414
$ws.Range("M102").Value = -611294.4399999999
$ws.Range("N102").Value = -28577658
$ws.Range("H132").Value = 2085.4285
$ws.Range("I132").Value = 2085.4285
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6256.2855
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3726.2855
$ws.Range("N132").Value = $null
$ws.Range("H136").Value = 3143.75
$ws.Range("I136").Value = 3143.75
$ws.Range("K136").Value = 9431.25
$ws.Range("M136").Value = -6881.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1506.5
$ws.Range("I80").Value = 2006
$ws.Range("J80").Value = 1007
$ws.Range("K80").Value = 2006
$ws.Range("L80").Value = 1007
$ws.Range("M80").Value = -1008
$ws.Range("N80").Value = -3003
$ws.Range("H83").Value = 1506.5
$ws.Range("I83").Value = 2006
$ws.Range("J83").Value = 1007
$ws.Range("K83").Value = 10030
$ws.Range("L83").Value = 5035
$ws.Range("M83").Value = -5038
$ws.Range("N83").Value = -15019
$ws.Range("H86").Value = 2270.2856
$ws.Range("I86").Value = 2148.8333
$ws.Range("K86").Value = 2148.8333
$ws.Range("M86").Value = -1025.8333
$ws.Range("H89").Value = 2270.2856
$ws.Range("I89").Value = 2148.8333
$ws.Range("K89").Value = 10744.1665
$ws.Range("M89").Value = -5128.166499999999
$ws.Range("H94").Value = 74736.664
$ws.Range("I94").Value = 92929.164
$ws.Range("K94").Value = 92929.164
$ws.Range("M94").Value = -92478.164
$ws.Range("H99").Value = 4796.85
$ws.Range("I99").Value = 4396.3335
$ws.Range("K99").Value = 4396.3335
$ws.Range("M99").Value = -2898.3335

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 19049.875
$ws.Range("J38").Value = 18779.8
$ws.Range("L38").Value = 18779.8
$ws.Range("N38").Value = -19533.8
$ws.Range("H46").Value = 19049.875
$ws.Range("J46").Value = 18779.8
$ws.Range("L46").Value = 18779.8
$ws.Range("N46").Value = -19201.8
$ws.Range("H107").Value = 671.3333
$ws.Range("I107").Value = 575.8
$ws.Range("J107").Value = 910.1667
$ws.Range("K107").Value = 575.8
$ws.Range("L107").Value = 910.1667
$ws.Range("M107").Value = 1344.2
$ws.Range("N107").Value = -4750.1667

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 19.40909
$ws.Range("I10").Value = 19.40909
$ws.Range("K10").Value = 58.22727
$ws.Range("M10").Value = 80.77273
$ws.Range("H12").Value = 29.444445
$ws.Range("I12").Value = 1.875
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 5.625
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = 167.375
$ws.Range("N12").Value = -1096
$ws.Range("H14").Value = 1792.4
$ws.Range("I14").Value = 1792.4
$ws.Range("K14").Value = 5377.200000000001
$ws.Range("M14").Value = -5204.200000000001
$ws.Range("H68").Value = 2335.8462
$ws.Range("I68").Value = 2168
$ws.Range("K68").Value = 6504
$ws.Range("M68").Value = -5693
$ws.Range("H71").Value = 2335.8462
$ws.Range("I71").Value = 2168
$ws.Range("K71").Value = 19512
$ws.Range("M71").Value = -15456
$ws.Range("H113").Value = 906.6
$ws.Range("J113").Value = 330
$ws.Range("L113").Value = 990
$ws.Range("N113").Value = -5330

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 99999.664
$ws.Range("J119").Value = 99999.664
$ws.Range("L119").Value = 99999.664
$ws.Range("N119").Value = -109675.664
$ws.Range("H132").Value = 1659.4
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1498.75
$ws.Range("I16").Value = 1498.75
$ws.Range("K16").Value = 1498.75
$ws.Range("M16").Value = -1328.75
$ws.Range("H43").Value = 8000
$ws.Range("J43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("N43").Value = -8386
$ws.Range("H70").Value = 29081.5
$ws.Range("J70").Value = 29081.5
$ws.Range("L70").Value = 29081.5
$ws.Range("N70").Value = -29621.5
$ws.Range("H73").Value = 29081.5
$ws.Range("J73").Value = 29081.5
$ws.Range("L73").Value = 29081.5
$ws.Range("N73").Value = -29621.5
$ws.Range("H122").Value = 3367.4443
$ws.Range("I122").Value = 2960.4
$ws.Range("J122").Value = 3876.25
$ws.Range("K122").Value = 8881.200000000001
$ws.Range("L122").Value = 11628.75
$ws.Range("M122").Value = -6431.200000000001
$ws.Range("N122").Value = -16528.75
$ws.Range("H132").Value = 1841.7
$ws.Range("I132").Value = 1133
$ws.Range("K132").Value = 3399
$ws.Range("M132").Value = -869
$ws.Range("H136").Value = 2998
$ws.Range("I136").Value = 2998
$ws.Range("K136").Value = 8994
$ws.Range("M136").Value = -6444

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 212.25
$ws.Range("I107").Value = 212.25
$ws.Range("K107").Value = 636.75
$ws.Range("M107").Value = 1283.25
